$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row appended at row 8
$ws.Range("A8").Value = 131180100
$ws.Range("B8").Value = 91804
$ws.Range("D8").Value = "NT"
$ws.Range("E8").Value = 1108
$ws.Range("F8").Value = "Harticka"
$ws.Range("G8").Value = "Pelloporus leporinus"
$ws.Range("H8").Value = "(Fr.) Krieglst."
$ws.Range("P8").Value = "Acksjön, Jmt"
$ws.Range("Q8").Value = 473246
$ws.Range("R8").Value = 7012290
$ws.Range("S8").Value = 100
$ws.Range("T8").Value = "Jämtland"
$ws.Range("U8").Value = "Krokom"
$ws.Range("V8").Value = "Jämtland"
$ws.Range("W8").Value = "Rödön"
# Dates are stored as plain text (not Excel date serials), keep them textual
$ws.Range("Y8").Value = "'2010-06-13"
$ws.Range("AA8").Value = "'2010-06-13"
$ws.Range("AD8").Value = $false
$ws.Range("AE8").Value = $false
$ws.Range("AG8").Value = $false
$ws.Range("AW8").Value = "Pontus Wallén"
$ws.Range("AX8").Value = "Pontus Wallén"
$ws.Range("AY8").Value = "LstZ naturvärdesinventeringar mellan åren xxxx-xxxx"
